# "data model of structures"
# Rebuild the "structures" sheet with the new INSPER_* data-model column
# layout, and rename the header of the "sections" sheet to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "structures": replace the 4-column layout with the new 20-column
# (A:T) data-model layout.
# ---------------------------------------------------------------------------
$structures = $wb.Worksheets.Item("structures")

$headers = @(
    "INSPER_ID_PRE",
    "BUSINESS_ID_PRE",
    "TYPE_OF_PARTICIPATION_CD",
    "TYPE_OF_INSURED_PERIOD_CD",
    "ACTIVE_FLAG_CD",
    "INSPER_EFFECTIVE_DATE",
    "INSPER_EXPIRY_DATE",
    "REPROG_ID_PRE",
    "BUSINESS_TITLE",
    "INSPER_LAYER_NO",
    "INSPER_MAIN_CURRENCY_CD",
    "INSPER_UW_YEAR",
    "INSPER_CONTRACT_ORDER",
    "INSPER_CONTRACT_FORM_CD_SLAV",
    "INSPER_CONTRACT_LODRA_CD_SLAV",
    "INSPER_CONTRACT_COVERAGE_CD_SLAV",
    "INSPER_CLAIM_BASIS_CD",
    "INSPER_LODRA_CD_SLAV",
    "INSPER_LOD_TO_RA_DATE_SLAV",
    "INSPER_COMMENT"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $structures.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# The original header formatting (bold font, border, centered) only covered
# columns A:D; extend it across the new header cells E1:T1 to match.
$structures.Range("A1").Copy()
$structures.Range("E1:T1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Clear out the old data columns beyond the header row before rewriting,
# in case any stray values remain outside the new layout.
$structures.Range("A2:D8").ClearContents()

$names = @("QS_1", "XOL_1", "XOL_2", "XOL_3", "XOL_4", "XOL_5", "XOL_6")
$types = @("quota_share", "excess_of_loss", "excess_of_loss", "excess_of_loss", "excess_of_loss", "excess_of_loss", "excess_of_loss")

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2

    $structures.Cells.Item($row, 1).Value = $i + 1                 # INSPER_ID_PRE
    $structures.Cells.Item($row, 2).Value = ""                     # BUSINESS_ID_PRE
    $structures.Cells.Item($row, 3).Value = $types[$i]             # TYPE_OF_PARTICIPATION_CD
    $structures.Cells.Item($row, 4).Value = ""                     # TYPE_OF_INSURED_PERIOD_CD
    $structures.Cells.Item($row, 5).Value = $true                  # ACTIVE_FLAG_CD
    $structures.Cells.Item($row, 6).Value = ""                     # INSPER_EFFECTIVE_DATE
    $structures.Cells.Item($row, 7).Value = ""                     # INSPER_EXPIRY_DATE
    $structures.Cells.Item($row, 8).Value = 1                      # REPROG_ID_PRE
    $structures.Cells.Item($row, 9).Value = $names[$i]              # BUSINESS_TITLE
    $structures.Cells.Item($row, 10).Value = ""                    # INSPER_LAYER_NO
    $structures.Cells.Item($row, 11).Value = ""                    # INSPER_MAIN_CURRENCY_CD
    $structures.Cells.Item($row, 12).Value = ""                    # INSPER_UW_YEAR
    $structures.Cells.Item($row, 13).Value = $i                    # INSPER_CONTRACT_ORDER
    $structures.Cells.Item($row, 14).Value = ""                    # INSPER_CONTRACT_FORM_CD_SLAV
    $structures.Cells.Item($row, 15).Value = ""                    # INSPER_CONTRACT_LODRA_CD_SLAV
    $structures.Cells.Item($row, 16).Value = ""                    # INSPER_CONTRACT_COVERAGE_CD_SLAV
    $structures.Cells.Item($row, 17).Value = "risk_attaching"       # INSPER_CLAIM_BASIS_CD
    $structures.Cells.Item($row, 18).Value = ""                    # INSPER_LODRA_CD_SLAV
    $structures.Cells.Item($row, 19).Value = ""                    # INSPER_LOD_TO_RA_DATE_SLAV
    $structures.Cells.Item($row, 20).Value = ""                    # INSPER_COMMENT
}

# ---------------------------------------------------------------------------
# Sheet "sections": rename the "structure_name" header to "BUSINESS_TITLE".
# ---------------------------------------------------------------------------
$sections = $wb.Worksheets.Item("sections")
$sections.Range("A1").Value = "BUSINESS_TITLE"
